# Shift the rolling "blog" slots in row 8 and append the new ser:129 post.
#
# Row 8 carries three "blog" widget cells whose text is:
#   C8 -> ser: 128
#   E8 -> ser: 127
#   I8 -> ser: 125
#
# The new post (ser: 129) pushes the window forward: the oldest entry
# (ser: 125) drops out, and every remaining entry shifts toward the
# "oldest" slot (I8), freeing C8 (the "newest" slot) for the new post.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 127"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 128"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 129"
